$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.477.52"
$ws.Range("E2").Value = "  +0.31%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.598.09"
$ws.Range("E3").Value = "  +0.84%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "516.04"
$ws.Range("E5").Value = "  +2.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.78"
$ws.Range("E6").Value = "  -0.20%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("E8").Value = "  +2.75%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.63"
$ws.Range("E9").Value = "  +0.21%  "
$ws.Range("E10").Value = "  +0.68%  "
$ws.Range("E11").Value = "  +0.92%  "
$ws.Range("E12").Value = "  +1.68%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.054.11"
$ws.Range("E13").Value = "  +1.18%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "60.538.35"
$ws.Range("E14").Value = "  +0.41%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.63"
$ws.Range("E15").Value = "  +0.03%  "
$ws.Range("E16").Value = "  +0.77%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.607.65"
$ws.Range("E17").Value = "  +0.72%  "
$ws.Range("E18").Value = "  -1.06%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "358.17"
$ws.Range("E19").Value = "  +3.49%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.58"
$ws.Range("E20").Value = "  +2.36%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.20"
$ws.Range("E21").Value = "  +2.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.998"
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "61.03"
$ws.Range("E23").Value = "  +1.43%  "
$ws.Range("E24").Value = "  +1.86%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.716.71"
$ws.Range("E25").Value = "  +1.34%  "
$ws.Range("E26").Value = "  +0.35%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.996"
$ws.Range("E27").Value = "  +0.40%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0838"
$ws.Range("E28").Value = "  -1.37%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.32"
$ws.Range("E29").Value = "  -2.10%  "
$ws.Range("E30").Value = "  +0.14%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "19.45"
$ws.Range("E31").Value = "  +1.27%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.03"
$ws.Range("E32").Value = "  +5.30%  "
$ws.Range("E33").Value = "  +2.18%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "148.39"
$ws.Range("E34").Value = "  -4.73%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.946"
$ws.Range("E35").Value = "  +9.87%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.01"
$ws.Range("E36").Value = "  +0.49%  "
$ws.Range("E37").Value = "  -0.42%  "
$ws.Range("E38").Value = "  +0.74%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.32"
$ws.Range("E39").Value = "  +2.05%  "
$ws.Range("E40").Value = "  +0.35%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.837"
$ws.Range("E41").Value = "  -0.96%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "288.23"
$ws.Range("E42").Value = "  -3.88%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.102"
$ws.Range("E43").Value = "  +1.69%  "
$ws.Range("E44").Value = "  +0.49%  "
$ws.Range("E45").Value = "  -1.42%  "
$ws.Range("E46").Value = "  +0.13%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.57"
$ws.Range("E47").Value = "  -0.99%  "
$ws.Range("E48").Value = "  +0.25%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0235"
$ws.Range("E49").Value = "  +0.83%  "
$ws.Range("E50").Value = "  +0.23%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.979.87"
$ws.Range("E51").Value = "  -2.41%  "
